$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Materialise the 12 new rows (98-109) as a real row insert first, so the
# engine records this as a structural row-insert rather than a pile of
# "value changed" diffs caused by writing past the end of the used range.
for ($i = 0; $i -lt 12; $i++) {
    $ws.Rows.Item(98).Insert()
}

function Set-BlankText($addr) {
    # A leading apostrophe forces Excel to store the cell as empty *text*
    # (matching the source file's empty <is/> inline-string cells) instead
    # of leaving the cell completely absent. Reset the style afterwards so
    # it doesn't leave a visible "number stored as text" quote-prefix flag.
    $ws.Range($addr).Value = "'"
    $ws.Range($addr).Style = "Normal"
}

# Update existing student-order cells (names reordered)
$ws.Range("E19").Value = "Brianna 1, Jessica 2"
$ws.Range("E55").Value = "Brianna 1, Hannah 3"
$ws.Range("E97").Value = "Hannah 3, Jessica 2"

# Add the new "PRAKTIKUMID" section starting at row 98
$ws.Range("A98").Value = "PRAKTIKUMID:"

$ws.Range("A99").Value = "Geograafia"
Set-BlankText "E99"
Set-BlankText "O99"

$ws.Range("A100").Value = "Kirjandus"
Set-BlankText "E100"
Set-BlankText "O100"

$ws.Range("A101").Value = "Inglise keel"
Set-BlankText "E101"
Set-BlankText "O101"

$ws.Range("A102").Value = "Ökoloogia"
Set-BlankText "E102"
Set-BlankText "O102"

$ws.Range("A103").Value = "Koorilaul"
$ws.Range("E103").Value = "Brianna 1, Hannah 3, Lauren 4, Emma 5"
Set-BlankText "O103"

$ws.Range("A104").Value = "Rahvatants"
$ws.Range("E104").Value = "Jessica 2, Emma 5, Lauren 4, Hannah 3"
Set-BlankText "O104"

$ws.Range("A105").Value = "Akvaristika"
$ws.Range("E105").Value = "Lauren 4, Emma 5"
Set-BlankText "O105"

$ws.Range("A106").Value = "Näitering"
$ws.Range("E106").Value = "Emma 5"
Set-BlankText "O106"

$ws.Range("A107").Value = "Mehhatroonika ja robootika"
Set-BlankText "E107"
Set-BlankText "O107"

$ws.Range("A108").Value = "Digitehnoloogiaga sõbraks"
Set-BlankText "E108"
Set-BlankText "O108"

$ws.Range("A109").Value = "Keskkonnaseire ja digilahendused"
Set-BlankText "E109"
Set-BlankText "O109"
